$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.460.45'
$ws.Range('E2').Value = '  +2.74%  '
$ws.Range('D3').Value = '2.358.36'
$ws.Range('E3').Value = '  +6.12%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.642'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.95%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('E9').Value = '  +6.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.18'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0939'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.87'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('E13').Value = '  +3.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.41'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +8.80%  '
$ws.Range('D16').Value = '2.714.84'
$ws.Range('E16').Value = '  +6.20%  '
$ws.Range('D17').Value = '2.419.79'
$ws.Range('E17').Value = '  +7.95%  '
$ws.Range('D18').Value = '43.449.50'
$ws.Range('E18').Value = '  +2.66%  '
$ws.Range('E19').Value = '  +3.62%  '
$ws.Range('E20').Value = '  -1.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '75.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '257.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +12.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('E26').Value = '  +3.54%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.55%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.10%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.59'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.55%  '
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '173.27'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0932'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.04'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.09%  '
$ws.Range('E35').Value = '  +5.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.98'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.15'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.12%  '
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.80'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +15.70%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +14.44%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.233'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('B44').Value = 'Celestia'
$ws.Range('C44').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.77'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.32%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('E46').Value = '  +3.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +11.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '111.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.70%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  +2.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.473'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.73%  '
